# feat: add 2022-Q3 data
#
# This script:
#  1. Inserts a brand new worksheet named "2022-Q3" right before the
#     existing "2022-Q1" sheet (so the sheet order becomes:
#     总计, 2022-Q3, 2022-Q1, 2021-Q4, 2021-Q3).
#  2. Populates the new sheet with the 2022-Q3 fund holding data.
#  3. Inserts a new row into the "总计" (totals) sheet for the 2022-Q3
#     quarter, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlCenter = -4108
$xlTop = -4160

# ---------------------------------------------------------------------
# Step 1: structural changes (must happen before we grab any other
# sheet references, since worksheet handles resolve by position and
# become stale once sheets are inserted/reordered).
# ---------------------------------------------------------------------

$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Add($q1Sheet)
$newSheet.Name = "2022-Q3"

$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# ---------------------------------------------------------------------
# Step 2: fill in the new "2022-Q3" sheet
# ---------------------------------------------------------------------

$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# Copy header (row1) + A-column (index) cell formatting from the
# "2022-Q1" sheet so the new sheet matches the look of its siblings.
$q1Sheet.Range("A1:H1").Copy()
$q3Sheet.Range("A1:H1").PasteSpecial($xlPasteFormats)
$q1Sheet.Range("A2:A5").Copy()
$q3Sheet.Range("A2:A9").PasteSpecial($xlPasteFormats)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Columns B (fund code), D, E, F, G hold numbers-as-text in the source
# data (e.g. leading zeros, fixed decimal points) so force them to
# text format before assigning, otherwise Excel will coerce them to
# numbers and mangle the representation (e.g. "002567" -> 2567).
$q3Sheet.Range("B2:B9").NumberFormat = "@"
$q3Sheet.Range("D2:G9").NumberFormat = "@"

$rows = @(
    @("001048", "富国新兴产业股票A",                 "17.84", "86.62", "4.58", "0.8171", 6),
    @("501077", "富国创新企业灵活配置混合（LOF）A",   "8.77",  "85.59", "7.47", "0.6551", 2),
    @("015686", "富国新兴产业股票C",                 "8.16",  "86.62", "4.58", "0.3737", 6),
    @("014611", "富国核心科技12个月持有期混合A",      "6.67",  "90.10", "2.16", "0.1441", 10),
    @("015133", "华安鼎安优选一年持有混合A",          "1.32",  "28.58", "1.05", "0.0139", 5),
    @("014612", "富国核心科技12个月持有期混合C",      "0.62",  "90.10", "2.16", "0.0134", 10),
    @("015134", "华安鼎安优选一年持有混合C",          "0.65",  "28.58", "1.05", "0.0068", 5),
    @("015849", "富国创新企业灵活配置混合（LOF）C",   "0.00",  "85.59", "7.47", $null,    2)
)

$r = 2
foreach ($row in $rows) {
    $q3Sheet.Cells.Item($r, 1).Value = $r - 2
    $q3Sheet.Cells.Item($r, 2).Value = $row[0]
    $q3Sheet.Cells.Item($r, 3).Value = $row[1]
    $q3Sheet.Cells.Item($r, 4).Value = $row[2]
    $q3Sheet.Cells.Item($r, 5).Value = $row[3]
    $q3Sheet.Cells.Item($r, 6).Value = $row[4]
    if ($null -eq $row[5]) {
        $q3Sheet.Cells.Item($r, 7).NumberFormat = "General"
        $q3Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q3Sheet.Cells.Item($r, 7).Value = $row[5]
    }
    $q3Sheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: update "总计" for the new 2022-Q3 quarter.
#
# The per-quarter rows are ordered most-recent-first; column A is just
# a 0-based position counter that does NOT change. Adding the new
# 2022-Q3 quarter therefore cascades the B/C/D (label/count/value)
# content of every existing row down by one, and a new row is appended
# for what used to be the oldest visible quarter (2021-Q3).
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 2.02

$totalSheet.Cells.Item(3, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(3, 3).Value = 4
$totalSheet.Cells.Item(3, 4).Value = 0.02

$totalSheet.Cells.Item(4, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(4, 3).Value = 4
$totalSheet.Cells.Item(4, 4).Value = 0.1

# New row 5 (A5 keeps the same column-A formatting as A2:A4).
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial($xlPasteFormats)
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 4
$totalSheet.Cells.Item(5, 4).Value = 0.89
